$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 20

# Capture the existing data before overwriting anything, since columns will shift:
#   old A (rows 2-20) -> category names, moves to new B
#   old B (rows 2-20, header "PercActivations") -> moves to new C
#   old C (rows 2-20, header "PercSegmentAreas") -> moves to new D
$oldHeaderB = $ws.Cells.Item(1, 2).Value2
$oldHeaderC = $ws.Cells.Item(1, 3).Value2

$names = @()
$valsB = @()
$valsC = @()
for ($r = 2; $r -le $lastRow; $r++) {
    $names += , $ws.Cells.Item($r, 1).Value2
    $valsB += , $ws.Cells.Item($r, 2).Value2
    $valsC += , $ws.Cells.Item($r, 3).Value2
}

# Propagate the header cell formatting (bold, bordered, centered style used by row 1)
# onto the new header cells B1 and D1 before they receive their new content; C1 already
# carries this formatting, so we reuse it as the source to avoid creating duplicate
# font/style entries.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("B1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Write new column D (was old C: "PercSegmentAreas")
$ws.Cells.Item(1, 4).Value = $oldHeaderC
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 4).Value = $valsC[$r - 2]
}

# Write new column C (was old B: "PercActivations")
$ws.Cells.Item(1, 3).Value = $oldHeaderB
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = $valsB[$r - 2]
}

# Write new column B: header "segments", data = category names (no special cell style)
$ws.Cells.Item(1, 2).Value = "segments"
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $names[$r - 2]
}

# Rewrite column A: data becomes 0-based numeric index (keeps its existing header style)
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Reset selection/active cell back to A1.
$ws.Range("A1").Select() | Out-Null
